# Applies the Dec 16 2022 symbol-list refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the target range to remain a text cell (matches the workbook's
    # original inline-string cells) instead of being auto-coerced to a number,
    # then drop back to the default "Normal" style so no stray formatting is
    # introduced (NumberFormat=@ would otherwise leave a Text number format).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Column G ("Hora") bumps from 16 to 17 for every data row (rows 2-51).
Set-TextValue $ws.Range("G2:G51") "17"

# Updated coin rows/prices/volumes (rankings reshuffled + new quotes).
Set-TextValue $ws.Range("D2") "242.99"
Set-TextValue $ws.Range("D3") "23.23"
Set-TextValue $ws.Range("D4") "5.731"
Set-TextValue $ws.Range("D5") "0.05815"
Set-TextValue $ws.Range("D6") "3.412"
Set-TextValue $ws.Range("D7") "6.486"
Set-TextValue $ws.Range("D8") "1.319"
Set-TextValue $ws.Range("D9") "0.8001"
Set-TextValue $ws.Range("B10") "WazirX"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1467"
Set-TextValue $ws.Range("E10") "9WazirXWRX"
Set-TextValue $ws.Range("B11") "MandalaExchangeToken"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.07623"
Set-TextValue $ws.Range("E11") "10MandalaExchangeTokenMDX"
Set-TextValue $ws.Range("B12") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D12") "0.03261"
Set-TextValue $ws.Range("E12") "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03013"
Set-TextValue $ws.Range("E13") "12BitrueCoinBTR"
Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09217"
Set-TextValue $ws.Range("E14") "13BitMartTokenBMX"
Set-TextValue $ws.Range("B15") "BitForexToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001663"
Set-TextValue $ws.Range("E15") "14BitForexTokenBF"
Set-TextValue $ws.Range("D16") "3.434"
Set-TextValue $ws.Range("B17") "CoinExToken"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D17") "0.04761"
Set-TextValue $ws.Range("E17") "16CoinExTokenCET"
Set-TextValue $ws.Range("B18") "One"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D18") "0.0006003"
Set-TextValue $ws.Range("E18") "17OneONE"
Set-TextValue $ws.Range("D19") "0.006264"
Set-TextValue $ws.Range("D20") "0.005457"
Set-TextValue $ws.Range("D21") "0.001072"
Set-TextValue $ws.Range("D23") "3.723"
Set-TextValue $ws.Range("D24") "2.210"
Set-TextValue $ws.Range("D25") "0.3333"
Set-TextValue $ws.Range("D40") "0.04306"
Set-TextValue $ws.Range("D41") "0.007038"
Set-TextValue $ws.Range("D42") "0.1050"
Set-TextValue $ws.Range("D43") "0.003408"
Set-TextValue $ws.Range("D44") "0.008618"
Set-TextValue $ws.Range("D45") "0.002471"
Set-TextValue $ws.Range("D46") "0.00005743"
Set-TextValue $ws.Range("D48") "0.7869"
Set-TextValue $ws.Range("D49") "0.1083"

Write-Host "Applied symbol list update (106 cells changed)"
